# Apply cryptos list update (Fri Oct 25 23:40:51 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.359.32"
$ws.Range("E2").Value = "  -2.58%  "

$ws.Range("D3").Value = "'2.415.58"
$ws.Range("E3").Value = "  -4.63%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'571.26"
$ws.Range("E5").Value = "  -3.78%  "

$ws.Range("D6").Value = "'164.04"
$ws.Range("E6").Value = "  -7.35%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.500"
$ws.Range("E8").Value = "  -5.98%  "

$ws.Range("D9").Value = "'2.413.54"
$ws.Range("E9").Value = "  -4.63%  "

$ws.Range("E10").Value = "  -8.29%  "

$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("D12").Value = "'0.323"
$ws.Range("E12").Value = "  -6.58%  "

$ws.Range("D13").Value = "'4.74"
$ws.Range("E13").Value = "  -7.47%  "

$ws.Range("D14").Value = "'24.64"
$ws.Range("E14").Value = "  -8.16%  "

$ws.Range("D15").Value = "'66.023.20"
$ws.Range("E15").Value = "  -3.11%  "

$ws.Range("E16").Value = "  -8.33%  "

$ws.Range("D17").Value = "'0.0000164"
$ws.Range("E17").Value = "  -8.76%  "

$ws.Range("D18").Value = "'2.379.54"
$ws.Range("E18").Value = "  -6.16%  "

$ws.Range("D19").Value = "'11.02"
$ws.Range("E19").Value = "  -4.58%  "

$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  -7.76%  "

$ws.Range("D21").Value = "'346.98"
$ws.Range("E21").Value = "  -5.52%  "

$ws.Range("D22").Value = "'3.95"
$ws.Range("E22").Value = "  -6.06%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'67.88"
$ws.Range("E24").Value = "  -4.25%  "

$ws.Range("D25").Value = "'4.14"
$ws.Range("E25").Value = "  -12.13%  "

$ws.Range("D26").Value = "'1.73"
$ws.Range("E26").Value = "  -10.65%  "

$ws.Range("D27").Value = "'8.89"
$ws.Range("E27").Value = "  -12.71%  "

$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").Value = "'2.520.56"
$ws.Range("E29").Value = "  -5.30%  "

$ws.Range("D30").Value = "'0.0₃0873"
$ws.Range("E30").Value = "  -12.40%  "

$ws.Range("D31").Value = "'7.67"
$ws.Range("E31").Value = "  -7.36%  "

$ws.Range("D32").Value = "'486.10"
$ws.Range("E32").Value = "  -10.20%  "

$ws.Range("D33").Value = "'1.75"
$ws.Range("E33").Value = "  -6.94%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "'1.19"
$ws.Range("E35").Value = "  -10.93%  "

$ws.Range("D36").Value = "'156.12"
$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.112"
$ws.Range("E37").Value = "  -13.60%  "

$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'18.54"
$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'18.07"
$ws.Range("E39").Value = "  -4.14%  "

$ws.Range("E40").Value = "  -9.58%  "

$ws.Range("D41").Value = "'1.63"
$ws.Range("E41").Value = "  -9.71%  "

$ws.Range("D42").Value = "'0.319"
$ws.Range("E42").Value = "  -10.15%  "

$ws.Range("D43").Value = "'4.58"
$ws.Range("E43").Value = "  -11.91%  "

$ws.Range("D44").Value = "'39.11"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  -10.08%  "

$ws.Range("D46").Value = "'136.82"
$ws.Range("E46").Value = "  -7.12%  "

$ws.Range("D47").Value = "'3.41"
$ws.Range("E47").Value = "  -8.53%  "

$ws.Range("D48").Value = "'0.503"
$ws.Range("E48").Value = "  -10.24%  "

$ws.Range("D49").Value = "'1.55"
$ws.Range("E49").Value = "  -8.80%  "

$ws.Range("D50").Value = "'0.0719"
$ws.Range("E50").Value = "  -4.95%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.573"
$ws.Range("E51").Value = "  -4.10%  "
